$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.871.52'
$ws.Range("E2").Value = '''  -0.25%  '

$ws.Range("D3").Value = '''3.146.10'
$ws.Range("E3").Value = '''  +0.24%  '

$ws.Range("E4").Value = '''  +0.03%  '

$ws.Range("D5").Value = '''592.90'
$ws.Range("E5").Value = '''  +0.17%  '

$ws.Range("D6").Value = '''145.48'
$ws.Range("E6").Value = '''  -1.27%  '

$ws.Range("E7").Value = '''  +0.05%  '

$ws.Range("D8").Value = '''3.137.57'
$ws.Range("E8").Value = '''  +0.16%  '

$ws.Range("E9").Value = '''  -0.71%  '

$ws.Range("E10").Value = '''  -0.46%  '

$ws.Range("D11").Value = '''5.88'
$ws.Range("E11").Value = '''  +1.96%  '

$ws.Range("D12").Value = '''0.461'
$ws.Range("E12").Value = '''  -1.71%  '

$ws.Range("E13").Value = '''  -2.96%  '

$ws.Range("D14").Value = '''37.22'
$ws.Range("E14").Value = '''  -0.52%  '

$ws.Range("D15").Value = '''3.665.28'
$ws.Range("E15").Value = '''  +0.34%  '

$ws.Range("E16").Value = '''  -1.38%  '

$ws.Range("D17").Value = '''7.34'
$ws.Range("E17").Value = '''  +2.17%  '

$ws.Range("D18").Value = '''3.141.41'
$ws.Range("E18").Value = '''  +0.15%  '

$ws.Range("D19").Value = '''63.732.29'
$ws.Range("E19").Value = '''  -0.24%  '

$ws.Range("D20").Value = '''468.60'
$ws.Range("E20").Value = '''  +0.40%  '

$ws.Range("D21").Value = '''14.38'
$ws.Range("E21").Value = '''  +0.00%  '

$ws.Range("D22").Value = '''0.734'
$ws.Range("E22").Value = '''  -0.32%  '

$ws.Range("E23").Value = '''  -0.66%  '

$ws.Range("B24").Value = 'Fetch.AI'
$ws.Range("C24").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D24").Value = '''2.35'
$ws.Range("E24").Value = '''  +7.54%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '''13.01'
$ws.Range("E25").Value = '''  -2.30%  '

$ws.Range("D26").Value = '''81.42'
$ws.Range("E26").Value = '''  -1.36%  '

$ws.Range("D28").Value = '''9.83'
$ws.Range("E28").Value = '''  +9.19%  '

$ws.Range("D29").Value = '''7.43'
$ws.Range("E29").Value = '''  +8.08%  '

$ws.Range("D30").Value = '''2.25'
$ws.Range("E30").Value = '''  +0.22%  '

$ws.Range("E31").Value = '''  -0.41%  '

$ws.Range("E32").Value = '''  +0.12%  '

$ws.Range("D33").Value = '''27.80'
$ws.Range("E33").Value = '''  +2.22%  '

$ws.Range("E34").Value = '''  +1.28%  '

$ws.Range("D35").Value = '''0.0₃0843'
$ws.Range("E35").Value = '''  -5.28%  '

$ws.Range("E36").Value = '''  +1.19%  '

$ws.Range("D37").Value = '''2.32'
$ws.Range("E37").Value = '''  -3.37%  '

$ws.Range("E38").Value = '''  +0.44%  '

$ws.Range("E39").Value = '''  -5.33%  '

$ws.Range("D40").Value = '''51.51'
$ws.Range("E40").Value = '''  +1.04%  '

$ws.Range("D41").Value = '''9.28'
$ws.Range("E41").Value = '''  +6.32%  '

$ws.Range("D42").Value = '''455.32'
$ws.Range("E42").Value = '''  -0.52%  '

$ws.Range("D43").Value = '''0.294'
$ws.Range("E43").Value = '''  +5.56%  '

$ws.Range("E44").Value = '''  -0.58%  '

$ws.Range("D45").Value = '''2.912.96'
$ws.Range("E45").Value = '''  +0.33%  '

$ws.Range("D46").Value = '''39.66'
$ws.Range("E46").Value = '''  +10.45%  '

$ws.Range("E47").Value = '''  -3.30%  '

$ws.Range("D48").Value = '''130.51'
$ws.Range("E48").Value = '''  +3.30%  '

$ws.Range("E50").Value = '''  +2.38%  '

$ws.Range("E51").Value = '''  -0.98%  '
